$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 58, pushing existing rows 58-109 down to 60-111.
$ws.Rows.Item(58).Insert()
$ws.Rows.Item(58).Insert()

# --- Populate new row 58 ---
$ws.Cells.Item(58,1).Value  = 1
$ws.Cells.Item(58,2).Value  = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(58,3).Value  = 'Arica y Parinacota'
$ws.Cells.Item(58,4).Value  = 44904
$ws.Cells.Item(58,5).Value  = 15
$ws.Cells.Item(58,6).Value  = 100114001
$ws.Cells.Item(58,7).Value  = 'Papa'
$ws.Cells.Item(58,8).Value  = 'Asterix'
$ws.Cells.Item(58,9).Value  = '1a (cosecha)'
$ws.Cells.Item(58,10).Value = 1000
$ws.Cells.Item(58,11).Value = 19000
$ws.Cells.Item(58,12).Value = 20000
$ws.Cells.Item(58,13).Value = 19500
$ws.Cells.Item(58,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(58,15).Value = 'Provincia de Melipilla'
$ws.Cells.Item(58,16).Value = 780
$ws.Cells.Item(58,17).Value = 25
$ws.Cells.Item(58,18).Value = 'Hortaliza'

# --- Populate new row 59 ---
$ws.Cells.Item(59,1).Value  = 1
$ws.Cells.Item(59,2).Value  = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(59,3).Value  = 'Arica y Parinacota'
$ws.Cells.Item(59,4).Value  = 44904
$ws.Cells.Item(59,5).Value  = 15
$ws.Cells.Item(59,6).Value  = 100114001
$ws.Cells.Item(59,7).Value  = 'Papa'
$ws.Cells.Item(59,8).Value  = 'Patagonia'
$ws.Cells.Item(59,9).Value  = '1a (cosecha)'
$ws.Cells.Item(59,10).Value = 1000
$ws.Cells.Item(59,11).Value = 16000
$ws.Cells.Item(59,12).Value = 17000
$ws.Cells.Item(59,13).Value = 16500
$ws.Cells.Item(59,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(59,15).Value = 'Región de Los Lagos'
$ws.Cells.Item(59,16).Value = 660
$ws.Cells.Item(59,17).Value = 25
$ws.Cells.Item(59,18).Value = 'Hortaliza'
